$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64, shifting the existing rows 64-75 down to 65-76.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new reading.
$ws.Cells.Item(64, 1).Value = 9
$ws.Cells.Item(64, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(64, 3).Value = "Metropolitana"
$ws.Cells.Item(64, 4).Value = 44769
$ws.Cells.Item(64, 5).Value = 13
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100102
$ws.Cells.Item(64, 8).Value = "Cítricos"
$ws.Cells.Item(64, 9).Value = 100102006
$ws.Cells.Item(64, 10).Value = "Pomelo"
$ws.Cells.Item(64, 11).Value = "Start Ruby"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 70
$ws.Cells.Item(64, 14).Value = 11000
$ws.Cells.Item(64, 15).Value = 12000
$ws.Cells.Item(64, 16).Value = 11429
$ws.Cells.Item(64, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(64, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 19).Value = 816
$ws.Cells.Item(64, 20).Value = 14

# Preserve the date-formatted style used by the rest of column D.
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
